# Fix language in Word docs
# Some of the styles in the reference docx files had US English as
# language. Remove the two superfluous empty paragraphs at the top of
# the document body and clear the hard-coded en-US language override on
# the "Lijst opsom.teken1" and "Lijstnummering1" list styles.

$d = $word.ActiveDocument

# --- 1. Drop the two extra blank paragraphs at the start of the body ---
# The body starts with three empty paragraphs; only the last one (the
# one that precedes the sectPr) should remain.
while ($d.Paragraphs.Count -gt 1) {
    $d.Paragraphs.Item(1).Range.Delete()
}

# --- 2. Clear the wrongly hard-coded en-US language on the list styles ---
$styleNames = @("Lijst opsom.teken1", "Lijstnummering1")
foreach ($styleName in $styleNames) {
    $style = $d.Styles($styleName)
    # wdLanguageNone (0) clears the forced English (US) proofing language
    # so the style no longer overrides the document's own language.
    $style.LanguageID = 0
}
